# Sample Project / Main.xlsx — update the "Rule" label of the last rule row
# (R40, cell B11 on the "Rules" sheet) to the text "1".
#
# The leading apostrophe forces Excel to store the value as literal text
# ("1") rather than auto-converting the numeric-looking entry into a
# number, matching the shared-string cell (t="s") produced by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"

$wb.Save()
